$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "M1"
$ws.Cells.Item(2, 2).Value = "Tnfsf14"
$ws.Cells.Item(2, 3).Value = "Tnfrsf14"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 2.786831
$ws.Cells.Item(2, 8).Value = 8.360493
$ws.Cells.Item(2, 9).Value = 0.1227769702371957
$ws.Cells.Item(2, 10).Value = 0.1227769702371957
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 3.853017333333334
$ws.Cells.Item(2, 14).Value = 11.559052
$ws.Cells.Item(2, 15).Value = 0.06452481780012836
$ws.Cells.Item(2, 16).Value = 0.06452481780012838
$ws.Cells.Item(2, 17).Value = 10.73770814807067
$ws.Cells.Item(2, 18).Value = 96.63937333263601
$ws.Cells.Item(2, 19).Value = 0.007922161634606835
$ws.Cells.Item(2, 20).Value = 0.007922161634606836

$ws.Cells.Item(3, 1).Value = "M1"
$ws.Cells.Item(3, 2).Value = "Tnfsf14"
$ws.Cells.Item(3, 3).Value = "Tnfrsf14"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 2.786831
$ws.Cells.Item(3, 8).Value = 8.360493
$ws.Cells.Item(3, 9).Value = 0.1227769702371957
$ws.Cells.Item(3, 10).Value = 0.1227769702371957
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 2.521954666666666
$ws.Cells.Item(3, 14).Value = 7.565863999999999
$ws.Cells.Item(3, 15).Value = 0.04223408598737598
$ws.Cells.Item(3, 16).Value = 0.04223408598737599
$ws.Cells.Item(3, 17).Value = 7.028261445661332
$ws.Cells.Item(3, 18).Value = 63.25435301095199
$ws.Cells.Item(3, 19).Value = 0.005185373118267223
$ws.Cells.Item(3, 20).Value = 0.005185373118267225

$ws.Cells.Item(4, 1).Value = "M1"
$ws.Cells.Item(4, 2).Value = "Tnfsf14"
$ws.Cells.Item(4, 3).Value = "Tnfrsf14"
$ws.Cells.Item(4, 4).Value = "M1"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 2.786831
$ws.Cells.Item(4, 8).Value = 8.360493
$ws.Cells.Item(4, 9).Value = 0.1227769702371957
$ws.Cells.Item(4, 10).Value = 0.1227769702371957
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 17.84381733333333
$ws.Cells.Item(4, 14).Value = 53.53145199999999
$ws.Cells.Item(4, 15).Value = 0.2988227050865691
$ws.Cells.Item(4, 16).Value = 0.2988227050865692
$ws.Cells.Item(4, 17).Value = 49.72770330287066
$ws.Cells.Item(4, 18).Value = 447.5493297258359
$ws.Cells.Item(4, 19).Value = 0.03668854636861199
$ws.Cells.Item(4, 20).Value = 0.036688546368612

$ws.Cells.Item(5, 1).Value = "M1"
$ws.Cells.Item(5, 2).Value = "Tnfsf14"
$ws.Cells.Item(5, 3).Value = "Tnfrsf14"
$ws.Cells.Item(5, 4).Value = "M2"
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 2.786831
$ws.Cells.Item(5, 8).Value = 8.360493
$ws.Cells.Item(5, 9).Value = 0.1227769702371957
$ws.Cells.Item(5, 10).Value = 0.1227769702371957
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 23.64775066666667
$ws.Cells.Item(5, 14).Value = 70.943252
$ws.Cells.Item(5, 15).Value = 0.3960186708606028
$ws.Cells.Item(5, 16).Value = 0.3960186708606028
$ws.Cells.Item(5, 17).Value = 65.90228463813733
$ws.Cells.Item(5, 18).Value = 593.120561743236
$ws.Cells.Item(5, 19).Value = 0.04862197256562602
$ws.Cells.Item(5, 20).Value = 0.04862197256562602

$ws.Cells.Item(6, 1).Value = "M1"
$ws.Cells.Item(6, 2).Value = "Tnfsf14"
$ws.Cells.Item(6, 3).Value = "Tnfrsf14"
$ws.Cells.Item(6, 4).Value = "Neutro"
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 0.6666666666666666
$ws.Cells.Item(6, 7).Value = 2.786831
$ws.Cells.Item(6, 8).Value = 8.360493
$ws.Cells.Item(6, 9).Value = 0.1227769702371957
$ws.Cells.Item(6, 10).Value = 0.1227769702371957
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 10.697805
$ws.Cells.Item(6, 14).Value = 32.093415
$ws.Cells.Item(6, 15).Value = 0.1791515217215829
$ws.Cells.Item(6, 16).Value = 0.1791515217215829
$ws.Cells.Item(6, 17).Value = 29.812974605955
$ws.Cells.Item(6, 18).Value = 268.316771453595
$ws.Cells.Item(6, 19).Value = 0.02199568105035911
$ws.Cells.Item(6, 20).Value = 0.02199568105035911

$ws.Cells.Item(7, 1).Value = "M1"
$ws.Cells.Item(7, 2).Value = "Tnfsf14"
$ws.Cells.Item(7, 3).Value = "Tnfrsf14"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 0.6666666666666666
$ws.Cells.Item(7, 7).Value = 2.786831
$ws.Cells.Item(7, 8).Value = 8.360493
$ws.Cells.Item(7, 9).Value = 0.1227769702371957
$ws.Cells.Item(7, 10).Value = 0.1227769702371957
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 1.149381666666667
$ws.Cells.Item(7, 14).Value = 3.448145
$ws.Cells.Item(7, 15).Value = 0.01924819854374075
$ws.Cells.Item(7, 16).Value = 0.01924819854374075
$ws.Cells.Item(7, 17).Value = 3.203132459498333
$ws.Cells.Item(7, 18).Value = 28.828192135485
$ws.Cells.Item(7, 19).Value = 0.002363235499724491
$ws.Cells.Item(7, 20).Value = 0.002363235499724492

$ws.Cells.Item(8, 1).Value = "M2"
$ws.Cells.Item(8, 2).Value = "Tnfsf14"
$ws.Cells.Item(8, 3).Value = "Tnfrsf14"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 2.556772666666667
$ws.Cells.Item(8, 8).Value = 7.670318
$ws.Cells.Item(8, 9).Value = 0.1126414919306584
$ws.Cells.Item(8, 10).Value = 0.1126414919306584
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 3.853017333333334
$ws.Cells.Item(8, 14).Value = 11.559052
$ws.Cells.Item(8, 15).Value = 0.06452481780012836
$ws.Cells.Item(8, 16).Value = 0.06452481780012838
$ws.Cells.Item(8, 17).Value = 9.851289402059557
$ws.Cells.Item(8, 18).Value = 88.661604618536
$ws.Cells.Item(8, 19).Value = 0.007268171743560365
$ws.Cells.Item(8, 20).Value = 0.007268171743560366

$ws.Cells.Item(9, 1).Value = "M2"
$ws.Cells.Item(9, 2).Value = "Tnfsf14"
$ws.Cells.Item(9, 3).Value = "Tnfrsf14"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 2.556772666666667
$ws.Cells.Item(9, 8).Value = 7.670318
$ws.Cells.Item(9, 9).Value = 0.1126414919306584
$ws.Cells.Item(9, 10).Value = 0.1126414919306584
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 2.521954666666666
$ws.Cells.Item(9, 14).Value = 7.565863999999999
$ws.Cells.Item(9, 15).Value = 0.04223408598737598
$ws.Cells.Item(9, 16).Value = 0.04223408598737599
$ws.Cells.Item(9, 17).Value = 6.448064758305778
$ws.Cells.Item(9, 18).Value = 58.032582824752
$ws.Cells.Item(9, 19).Value = 0.004757310455945745
$ws.Cells.Item(9, 20).Value = 0.004757310455945746

$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Tnfsf14"
$ws.Cells.Item(10, 3).Value = "Tnfrsf14"
$ws.Cells.Item(10, 4).Value = "M1"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 2.556772666666667
$ws.Cells.Item(10, 8).Value = 7.670318
$ws.Cells.Item(10, 9).Value = 0.1126414919306584
$ws.Cells.Item(10, 10).Value = 0.1126414919306584
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 17.84381733333333
$ws.Cells.Item(10, 14).Value = 53.53145199999999
$ws.Cells.Item(10, 15).Value = 0.2988227050865691
$ws.Cells.Item(10, 16).Value = 0.2988227050865692
$ws.Cells.Item(10, 17).Value = 45.62258442685955
$ws.Cells.Item(10, 18).Value = 410.603259841736
$ws.Cells.Item(10, 19).Value = 0.03365983532370629
$ws.Cells.Item(10, 20).Value = 0.0336598353237063

$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Tnfsf14"
$ws.Cells.Item(11, 3).Value = "Tnfrsf14"
$ws.Cells.Item(11, 4).Value = "M2"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 2.556772666666667
$ws.Cells.Item(11, 8).Value = 7.670318
$ws.Cells.Item(11, 9).Value = 0.1126414919306584
$ws.Cells.Item(11, 10).Value = 0.1126414919306584
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 23.64775066666667
$ws.Cells.Item(11, 14).Value = 70.943252
$ws.Cells.Item(11, 15).Value = 0.3960186708606028
$ws.Cells.Item(11, 16).Value = 0.3960186708606028
$ws.Cells.Item(11, 17).Value = 60.46192253268178
$ws.Cells.Item(11, 18).Value = 544.1573027941361
$ws.Cells.Item(11, 19).Value = 0.04460813391813467
$ws.Cells.Item(11, 20).Value = 0.04460813391813467

$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Tnfsf14"
$ws.Cells.Item(12, 3).Value = "Tnfrsf14"
$ws.Cells.Item(12, 4).Value = "Neutro"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 2.556772666666667
$ws.Cells.Item(12, 8).Value = 7.670318
$ws.Cells.Item(12, 9).Value = 0.1126414919306584
$ws.Cells.Item(12, 10).Value = 0.1126414919306584
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 10.697805
$ws.Cells.Item(12, 14).Value = 32.093415
$ws.Cells.Item(12, 15).Value = 0.1791515217215829
$ws.Cells.Item(12, 16).Value = 0.1791515217215829
$ws.Cells.Item(12, 17).Value = 27.35185541733
$ws.Cells.Item(12, 18).Value = 246.16669875597
$ws.Cells.Item(12, 19).Value = 0.02017989468836686
$ws.Cells.Item(12, 20).Value = 0.02017989468836686

$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Tnfsf14"
$ws.Cells.Item(13, 3).Value = "Tnfrsf14"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 2.556772666666667
$ws.Cells.Item(13, 8).Value = 7.670318
$ws.Cells.Item(13, 9).Value = 0.1126414919306584
$ws.Cells.Item(13, 10).Value = 0.1126414919306584
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 1.149381666666667
$ws.Cells.Item(13, 14).Value = 3.448145
$ws.Cells.Item(13, 15).Value = 0.01924819854374075
$ws.Cells.Item(13, 16).Value = 0.01924819854374075
$ws.Cells.Item(13, 17).Value = 2.938707628901112
$ws.Cells.Item(13, 18).Value = 26.44836866011
$ws.Cells.Item(13, 19).Value = 0.002168145800944485
$ws.Cells.Item(13, 20).Value = 0.002168145800944485

$ws.Cells.Item(14, 1).Value = "Neutro"
$ws.Cells.Item(14, 2).Value = "Tnfsf14"
$ws.Cells.Item(14, 3).Value = "Tnfrsf14"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 17.35471666666666
$ws.Cells.Item(14, 8).Value = 52.06415
$ws.Cells.Item(14, 9).Value = 0.7645815378321459
$ws.Cells.Item(14, 10).Value = 0.7645815378321459
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 3.853017333333334
$ws.Cells.Item(14, 14).Value = 11.559052
$ws.Cells.Item(14, 15).Value = 0.06452481780012836
$ws.Cells.Item(14, 16).Value = 0.06452481780012838
$ws.Cells.Item(14, 17).Value = 66.86802413175556
$ws.Cells.Item(14, 18).Value = 601.8122171858
$ws.Cells.Item(14, 19).Value = 0.04933448442196117
$ws.Cells.Item(14, 20).Value = 0.04933448442196118

$ws.Cells.Item(15, 1).Value = "Neutro"
$ws.Cells.Item(15, 2).Value = "Tnfsf14"
$ws.Cells.Item(15, 3).Value = "Tnfrsf14"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 17.35471666666666
$ws.Cells.Item(15, 8).Value = 52.06415
$ws.Cells.Item(15, 9).Value = 0.7645815378321459
$ws.Cells.Item(15, 10).Value = 0.7645815378321459
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 2.521954666666666
$ws.Cells.Item(15, 14).Value = 7.565863999999999
$ws.Cells.Item(15, 15).Value = 0.04223408598737598
$ws.Cells.Item(15, 16).Value = 0.04223408598737599
$ws.Cells.Item(15, 17).Value = 43.76780868617777
$ws.Cells.Item(15, 18).Value = 393.9102781756
$ws.Cells.Item(15, 19).Value = 0.03229140241316301
$ws.Cells.Item(15, 20).Value = 0.03229140241316302

$ws.Cells.Item(16, 1).Value = "Neutro"
$ws.Cells.Item(16, 2).Value = "Tnfsf14"
$ws.Cells.Item(16, 3).Value = "Tnfrsf14"
$ws.Cells.Item(16, 4).Value = "M1"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 17.35471666666666
$ws.Cells.Item(16, 8).Value = 52.06415
$ws.Cells.Item(16, 9).Value = 0.7645815378321459
$ws.Cells.Item(16, 10).Value = 0.7645815378321459
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 17.84381733333333
$ws.Cells.Item(16, 14).Value = 53.53145199999999
$ws.Cells.Item(16, 15).Value = 0.2988227050865691
$ws.Cells.Item(16, 16).Value = 0.2988227050865692
$ws.Cells.Item(16, 17).Value = 309.6743940717554
$ws.Cells.Item(16, 18).Value = 2787.0695466458
$ws.Cells.Item(16, 19).Value = 0.2284743233942508
$ws.Cells.Item(16, 20).Value = 0.2284743233942509

$ws.Cells.Item(17, 1).Value = "Neutro"
$ws.Cells.Item(17, 2).Value = "Tnfsf14"
$ws.Cells.Item(17, 3).Value = "Tnfrsf14"
$ws.Cells.Item(17, 4).Value = "M2"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 17.35471666666666
$ws.Cells.Item(17, 8).Value = 52.06415
$ws.Cells.Item(17, 9).Value = 0.7645815378321459
$ws.Cells.Item(17, 10).Value = 0.7645815378321459
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 23.64775066666667
$ws.Cells.Item(17, 14).Value = 70.943252
$ws.Cells.Item(17, 15).Value = 0.3960186708606028
$ws.Cells.Item(17, 16).Value = 0.3960186708606028
$ws.Cells.Item(17, 17).Value = 410.4000126239777
$ws.Cells.Item(17, 18).Value = 3693.6001136158
$ws.Cells.Item(17, 19).Value = 0.3027885643768421
$ws.Cells.Item(17, 20).Value = 0.3027885643768422

$ws.Cells.Item(18, 1).Value = "Neutro"
$ws.Cells.Item(18, 2).Value = "Tnfsf14"
$ws.Cells.Item(18, 3).Value = "Tnfrsf14"
$ws.Cells.Item(18, 4).Value = "Neutro"
$ws.Cells.Item(18, 5).Value = 3
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = 17.35471666666666
$ws.Cells.Item(18, 8).Value = 52.06415
$ws.Cells.Item(18, 9).Value = 0.7645815378321459
$ws.Cells.Item(18, 10).Value = 0.7645815378321459
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 12).Value = 1
$ws.Cells.Item(18, 13).Value = 10.697805
$ws.Cells.Item(18, 14).Value = 32.093415
$ws.Cells.Item(18, 15).Value = 0.1791515217215829
$ws.Cells.Item(18, 16).Value = 0.1791515217215829
$ws.Cells.Item(18, 17).Value = 185.65737473025
$ws.Cells.Item(18, 18).Value = 1670.91637257225
$ws.Cells.Item(18, 19).Value = 0.136975945982857
$ws.Cells.Item(18, 20).Value = 0.136975945982857

$ws.Cells.Item(19, 1).Value = "Neutro"
$ws.Cells.Item(19, 2).Value = "Tnfsf14"
$ws.Cells.Item(19, 3).Value = "Tnfrsf14"
$ws.Cells.Item(19, 4).Value = "sCs"
$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = 17.35471666666666
$ws.Cells.Item(19, 8).Value = 52.06415
$ws.Cells.Item(19, 9).Value = 0.7645815378321459
$ws.Cells.Item(19, 10).Value = 0.7645815378321459
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 12).Value = 1
$ws.Cells.Item(19, 13).Value = 1.149381666666667
$ws.Cells.Item(19, 14).Value = 3.448145
$ws.Cells.Item(19, 15).Value = 0.01924819854374075
$ws.Cells.Item(19, 16).Value = 0.01924819854374075
$ws.Cells.Item(19, 17).Value = 19.94719316686111
$ws.Cells.Item(19, 18).Value = 179.52473850175
$ws.Cells.Item(19, 19).Value = 0.01471681724307177
$ws.Cells.Item(19, 20).Value = 0.01471681724307178
